# SizingGuide.xlsx edit: "Changed is_nan() def and added extrapolation"
#
# 1. Archer sheet: change base cooler dims (A2/B2), simplify the
#    C24/C25/C26 sizing formulas (drop the $B$20 absolute ref / IF wrapper).
# 2. Bulldozer CSF sheet: double the fin-pitch constant B20 (5 -> 10),
#    which ripples through the sizing/extrapolation table.
# 3. Move the active selection/tab from Archer to Bulldozer CSF.

$wb = $excel.ActiveWorkbook

$archer = $wb.Worksheets.Item("Archer")
$bulldozer = $wb.Worksheets.Item("Bulldozer CSF")

# --- Archer: base heatsink fin dimensions ---
$archer.Range("A2").Value = 8
$archer.Range("B2").Value = 2

# --- Archer: is_nan()-style formula cleanup on the C24:C26 sizing cells ---
$archer.Range("C24").Formula = "=B20*B24"
$archer.Range("C25").Formula = "=B20*B25"
$archer.Range("C26").Formula = "=B20*B26"

# --- Bulldozer CSF: extrapolate fin pitch ---
$bulldozer.Range("B20").Value = 10

# --- View state: Archer loses the tab selection / moves its cell cursor ---
[void]$archer.Range("C5").Select()

# --- View state: Bulldozer CSF becomes the active sheet/tab with A5:C5 selected ---
$bulldozer.Activate()
[void]$bulldozer.Range("A5:C5").Select()
